$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.661.55'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '1.905.55'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '''239.32'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').Value = '''0.4739'
$ws.Range('E7').Value = '  -0.56%  '
$ws.Range('D8').Value = '''0.2863'
$ws.Range('E8').Value = '  +0.30%  '
$ws.Range('D9').Value = '''0.06666'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = '''19.74'
$ws.Range('E10').Value = '  +5.11%  '
$ws.Range('D11').Value = '''101.24'
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('D12').Value = '''0.07806'
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('D13').Value = '1.910.44'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').Value = '''5.181'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').Value = '''0.6802'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '''288.01'
$ws.Range('E16').Value = '  +11.38%  '
$ws.Range('D17').Value = '30.656.34'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '''0.000007511'
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.161.32'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').Value = '''12.74'
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('D22').Value = '''5.422'
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = '''1.002'
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '''6.275'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '''9.361'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''167.15'
$ws.Range('E26').Value = '  +2.72%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''19.40'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '''2.036'
$ws.Range('E28').Value = '  -1.72%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''1.384'
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '''0.09949'
$ws.Range('E30').Value = '  -1.36%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '''4.519'
$ws.Range('E31').Value = '  -1.53%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '''1.513'
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '''4.248'
$ws.Range('E33').Value = '  +1.00%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '''0.04761'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '''0.7247'
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '''1.111'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').Value = '''1.001'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').Value = '''2.722'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').Value = '''0.01910'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('D40').Value = '''6.803'
$ws.Range('E40').Value = '  +8.52%  '
$ws.Range('D41').Value = '''2.601'
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('D42').Value = '''74.20'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = '''1.993'
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').Value = '''0.8702'
$ws.Range('E44').Value = '  +1.09%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '''104.93'
$ws.Range('E45').Value = '  -1.61%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').Value = '''0.4278'
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('D47').Value = '''1.001'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('D48').Value = '1.004.98'
$ws.Range('E48').Value = '  -0.31%  '
$ws.Range('D49').Value = '''7.388'
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('D50').Value = '''9.228'
$ws.Range('E50').Value = '  +4.75%  '
$ws.Range('D51').Value = '''0.1183'
$ws.Range('E51').Value = '  -0.72%  '
